# Update cryptos list with refreshed prices / volumes (and a handful of
# rank swaps) matching the Dec 9 2023 GitHub Actions data refresh.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: row, Coin, Link, Price, Volume(1h)
$data = @(
    @(2,  "Bitcoin",                     "https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc",                    "44.079.47",  "  +0.56%  "),
    @(3,  "Ethereum",                    "https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth",                   "2.354.91",   "  -0.08%  "),
    @(4,  "TetherUSD",                   "https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt",                 "1.00",       "  +0.17%  "),
    @(5,  "XRP",                         "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp",                        "0.681",      "  +1.37%  "),
    @(6,  "BNB",                         "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb",                        "239.64",     "  +1.68%  "),
    @(7,  "Solana",                      "https://coinranking.com/coin/zNZHO_Sjf+solana-sol",                         "74.70",      "  +1.89%  "),
    @(8,  "USDC",                        "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc",                      "1.00",       "  +0.06%  "),
    @(9,  "Cardano",                     "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada",                    "0.598",      "  +6.81%  "),
    @(10, "Dogecoin",                    "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge",                  "0.101",      "  +2.19%  "),
    @(11, "OKB",                         "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb",                        "57.22",      "  +0.06%  "),
    @(12, "Avalanche",                   "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax",                     "32.44",      "  +15.47%  "),
    @(13, "TRON",                        "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx",                       "0.108",      "  +1.02%  "),
    @(14, "Polkadot",                    "https://coinranking.com/coin/25W7FG7om+polkadot-dot",                       "7.24",       "  +7.44%  "),
    @(15, "WrappedliquidstakedEther2.0", "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth",  "2.709.35",   "  +0.04%  "),
    @(16, "Chainlink",                   "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link",                 "16.55",      "  -0.94%  "),
    @(17, "Polygon",                     "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic",                   "0.903",      "  +1.78%  "),
    @(18, "WrappedEther",                "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth",              "2.358.22",   "  +1.64%  "),
    @(19, "WrappedBTC",                  "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc",                 "44.107.82",  "  +0.62%  "),
    @(20, "ShibaInu",                    "https://coinranking.com/coin/xz24e0BjL+shibainu-shib",                      "0.0000103",  "  +1.50%  "),
    @(21, "Uniswap",                     "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni",                        "6.73",       "  +5.01%  "),
    @(22, "Litecoin",                    "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc",                   "76.99",      "  -0.99%  "),
    @(23, "BitcoinCash",                 "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch",                "256.89",     "  +1.19%  "),
    @(24, "Dai",                         "https://coinranking.com/coin/MoTuySvg7+dai-dai",                            "1.00",       "  +0.07%  "),
    @(25, "ImmutableX",                  "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx",                     "1.88",       "  +18.26%  "),
    @(26, "WEMIXToken",                  "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix",                   "3.71",       "  -1.41%  "),
    @(27, "PancakeSwap",                 "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake",                   "2.50",       "  +0.41%  "),
    @(28, "Cosmos",                      "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom",                    "10.74",      "  +0.96%  "),
    @(29, "Toncoin",                     "https://coinranking.com/coin/67YlI0K1b+toncoin-ton",                        "2.24",       "  -2.05%  "),
    @(30, "EthereumClassic",             "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc",            "22.85",      "  +1.72%  "),
    @(31, "Monero",                      "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr",                     "175.07",     "  +1.44%  "),
    @(32, "Stellar",                     "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm",                    "0.138",      "  +3.41%  "),
    @(33, "Kaspa",                       "https://coinranking.com/coin/V8GxkwWow+kaspa-kas",                          "0.127",      "  -2.37%  "),
    @(34, "Hedera",                      "https://coinranking.com/coin/jad286TjB+hedera-hbar",                        "0.0761",     "  +5.78%  "),
    @(35, "Filecoin",                    "https://coinranking.com/coin/ymQub4fuB+filecoin-fil",                       "5.29",       "  +1.76%  "),
    @(36, "InternetComputer(DFINITY)",   "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp",        "5.37",       "  +3.60%  "),
    @(37, "RenderToken",                 "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr",               "3.74",       "  -1.46%  "),
    @(38, "LidoDAOToken",                "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo",                   "2.37",       "  -3.01%  "),
    @(39, "THORChain",                   "https://coinranking.com/coin/ybmU-kKU+thorchain-rune",                      "6.36",       "  -0.84%  "),
    @(40, "VeChain",                     "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet",                    "0.0282",     "  +4.62%  "),
    @(41, "InjectiveProtocol",           "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj",              "19.30",      "  +0.25%  "),
    @(42, "Cronos",                      "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro",                      "0.108",      "  +10.40%  "),
    @(43, "Algorand",                    "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo",                  "0.205",      "  +10.83%  "),
    @(44, "FraxShare",                   "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs",                      "9.06",       "  +1.32%  "),
    @(45, "BinanceUSD",                  "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd",                "1.00",       "  +0.12%  "),
    @(46, "FTXToken",                    "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt",                       "4.73",       "  +6.35%  "),
    @(47, "NEARProtocol",                "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near",                  "2.50",       "  +8.57%  "),
    @(48, "TrustWalletToken",            "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt",               "1.25",       "  +2.54%  "),
    @(49, "MultiversX",                  "https://coinranking.com/coin/omwkOTglq+multiversx-egld",                   "57.15",      "  +9.55%  "),
    @(50, "Aave",                        "https://coinranking.com/coin/ixgUfzmLR+aave-aave",                          "100.86",     "  +3.38%  "),
    @(51, "ARBITRUM",                    "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb",                       "1.17",       "  +0.90%  ")
)

foreach ($entry in $data) {
    $row   = $entry[0]
    $coin  = $entry[1]
    $link  = $entry[2]
    $price = $entry[3]
    $vol   = $entry[4]

    $ws.Cells.Item($row, 2).Value = $coin
    $ws.Cells.Item($row, 3).Value = $link

    # Price column holds values like "44.131.77" / "0.683" / "1.00" that are
    # meant to stay plain text (they are not valid numbers, or would lose
    # their formatting such as trailing zeros if Excel auto-converted them).
    # Prefixing with an apostrophe forces text entry; resetting the style
    # back to "Normal" afterwards avoids leaving a stray quote-prefix style
    # behind, since the source cells carry no explicit style.
    $priceCell = $ws.Cells.Item($row, 4)
    $priceCell.Value = "'" + $price
    $priceCell.Style = "Normal"

    $ws.Cells.Item($row, 5).Value = $vol
}
